# Apply the three wording corrections requested in Leanne's review.
#
# Note: Find.Execute's inline "Replace" argument runs the replacement text
# through Word's AutoFormat/AutoCorrect pipeline, which silently turns plain
# straight apostrophes (') into curly/smart quotes (U+2019). The source text
# uses plain apostrophes, so instead we locate each target range with
# Find.Execute (no replacement argument) and then assign Range.Text directly,
# which performs a literal text substitution without any autocorrect.

$d = $word.ActiveDocument

# 1. "...(replacing student number by your actual student number)." ->
#    "...(replacing studentnumber by your actual student number and q1 with q2 or q3 for each of the questions)."
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    'Write three separate scripts for each question naming them as follows: `ma1003-studentnumber-q1.py'' (replacing student number by your actual student number). Submission details are given at the end of this sheet.',
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0
)
if ($found1) {
    $rng1.Text = 'Write three separate scripts for each question naming them as follows: `ma1003-studentnumber-q1.py'' (replacing studentnumber by your actual student number and q1 with q2 or q3 for each of the questions). Submission details are given at the end of this sheet.'
}

# 2. "Save your 3 files to the directory xyz." ->
#    "Save your 3 files to the MA1003 class test folder in the shared drive."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    'Save your 3 files to the directory xyz. You will not be able to see your file when it is saved there but be sure to save it there (perhaps do this twice).',
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0
)
if ($found2) {
    $rng2.Text = 'Save your 3 files to the MA1003 class test folder in the shared drive. You will not be able to see your file when it is saved there but be sure to save it there (perhaps do this twice).'
}

# 3. Fix typo "MAT1003-classtest-c123456789" -> "MA1003-classtest-c123456789"
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "MA1003-classtest-yourstudentnumber. For example: 'MAT1003-classtest-c123456789'.",
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0
)
if ($found3) {
    $rng3.Text = "MA1003-classtest-yourstudentnumber. For example: 'MA1003-classtest-c123456789'."
}
